$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column Q: header + first two data rows' new property values
$ws.Range("Q2").Value = "veryLongPropertyName"
$ws.Range("Q3").Value = "abc"
$ws.Range("Q4").Value = "xyz"

# New block of rows 19-21 describing another data.Foo entry
$ws.Range("A19").Value = "data.Foo"
$ws.Range("A20").Value = "nr"
$ws.Range("B20").Value = "vlpn"
$ws.Range("A21").Value = 10
$ws.Range("B21").Value = "ABC"

# Move selection to B21 like the final author state
$null = $ws.Range("B21").Select()
